$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ELIC-300 / ULBRICH M (row 19): corrected grade-distribution percentages
$ws.Range("C19").Value = 2.0477
$ws.Range("D19").Value = "2.08%"
$ws.Range("E19").Value = "47.92%"
$ws.Range("F19").Value = "22.92%"
$ws.Range("G19").Value = "6.25%"
$ws.Range("H19").Value = "20.83%"

# ELIG-200 / PAYNE G (row 32): corrected grade-distribution percentages
$ws.Range("C32").Value = 3.0625
$ws.Range("D32").Value = "34.38%"
$ws.Range("F32").Value = "9.38%"
